$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so Excel
# does not auto-convert numeric-looking strings (e.g. "0.999")
# into real numbers. Style is reset back to Normal afterwards so
# no stray formatting is introduced on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '51.140.91'
$ws.Range("E2").Value = '  -0.63%  '

$ws.Range("D3").Value = '2.955.13'
$ws.Range("E3").Value = '  +0.85%  '

Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.11%  '

Set-TextValue $ws.Range("D5") '380.03'
$ws.Range("E5").Value = '  +0.72%  '

Set-TextValue $ws.Range("D6") '102.18'
$ws.Range("E6").Value = '  -1.87%  '

$ws.Range("E7").Value = '  -0.62%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +0.35%  '

Set-TextValue $ws.Range("D10") '36.53'
$ws.Range("E10").Value = '  -1.29%  '

$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = '3.418.72'
$ws.Range("E13").Value = '  +0.48%  '

Set-TextValue $ws.Range("D14") '18.04'
$ws.Range("E14").Value = '  -2.09%  '

Set-TextValue $ws.Range("D15") '7.42'
$ws.Range("E15").Value = '  +0.48%  '

$ws.Range("D16").Value = '2.945.28'
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("E17").Value = '  +4.26%  '

$ws.Range("D18").Value = '51.068.45'
$ws.Range("E18").Value = '  -0.79%  '

Set-TextValue $ws.Range("D19") '3.20'
$ws.Range("E19").Value = '  -5.43%  '

Set-TextValue $ws.Range("D20") '7.12'
$ws.Range("E20").Value = '  -2.91%  '

Set-TextValue $ws.Range("D21") '12.54'
$ws.Range("E21").Value = '  -3.72%  '

$ws.Range("E22").Value = '  +0.24%  '

Set-TextValue $ws.Range("D23") '68.47'
$ws.Range("E23").Value = '  +0.15%  '

Set-TextValue $ws.Range("D24") '261.72'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("E25").Value = '  +2.26%  '

Set-TextValue $ws.Range("D26") '8.33'
$ws.Range("E26").Value = '  +13.19%  '

$ws.Range("E27").Value = '  +2.84%  '

$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("E29").Value = '  -0.73%  '

$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("E31").Value = '  +9.38%  '

Set-TextValue $ws.Range("D32") '25.65'
$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("E33").Value = '  -0.39%  '

Set-TextValue $ws.Range("D34") '0.0461'
$ws.Range("E34").Value = '  +7.80%  '

$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D35") '33.91'
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("E36").Value = '  -1.88%  '

$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D37") '50.42'
$ws.Range("E37").Value = '  -3.77%  '

$ws.Range("E38").Value = '  -0.29%  '

$ws.Range("E39").Value = '  -1.22%  '

Set-TextValue $ws.Range("D40") '16.79'
$ws.Range("E40").Value = '  -1.83%  '

$ws.Range("E41").Value = '  -2.49%  '

$ws.Range("E42").Value = '  +0.93%  '

$ws.Range("E43").Value = '  -2.38%  '

Set-TextValue $ws.Range("D44") '121.47'
$ws.Range("E44").Value = '  -0.52%  '

Set-TextValue $ws.Range("D45") '21.15'
$ws.Range("E45").Value = '  -2.82%  '

Set-TextValue $ws.Range("D46") '2.06'
$ws.Range("E46").Value = '  -0.26%  '

$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("D49").Value = '2.005.70'
$ws.Range("E49").Value = '  -0.68%  '

Set-TextValue $ws.Range("D50") '3.23'
$ws.Range("E50").Value = '  +1.67%  '

Set-TextValue $ws.Range("D51") '0.0339'
$ws.Range("E51").Value = '  +5.48%  '
